# "Update hours for monday" - the Agenda sheet tracks hours per subject per
# weekday. Row 2 is Monday. Column E is "Python", column H is "Pre-Calculus".
# Both were updated from 0 hours to 1 hour worked on Monday.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("H2").Value = 1

# Reflect the new cell selection recorded in the saved workbook.
$ws.Range("H14").Select()
